$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = '{''(the'', ''journey)''}'
$ws.Cells.Item(5, 2).Value = '{''spritz'', ''hotelwe'', ''aperol'', ''usthe'', ''lemoncello'', ''spritzi'', ''dinnertimethe'', ''himthey'', ''lightsthe''}'
$ws.Cells.Item(7, 2).Value = '{''noneall'', ''excellentimmaculately'', ''alsothe'', ''viewbut'', ''standardtoiletries'', ''aminta''}'
$ws.Cells.Item(10, 2).Value = '{''â‚¬20'', ''areathe'', ''bottlethe''}'
$ws.Cells.Item(11, 2).Value = '{''conection'', ''217was'', ''greatmy''}'
$ws.Cells.Item(12, 2).Value = '{''aminta'', ''courtsey'', ''relaxthere''}'
$ws.Cells.Item(14, 2).Value = '{''aminta'', ''day/evening'', ''itâ€™s''}'
$ws.Cells.Item(15, 2).Value = '{''aminta'', ''breakfest''}'
$ws.Cells.Item(16, 2).Value = '{''all)'', ''(this'', ''iâ€™d'', ''car()'', ''itâ€™s'', ''late)overall'', ''(and''}'
$ws.Cells.Item(17, 2).Value = '{''veggie'', ''shuttle-service'', ''cleanthe'', ''somethingsfood''}'
$ws.Cells.Item(20, 2).Value = '{''sorrentoour'', ''worldand'', ''finallt'', ''wifeâ€™s'', ''travellerwe'', ''breathtakingwe'', ''caprialthough'', ''coastnaples'', ''maintenancedo''}'
$ws.Cells.Item(24, 2).Value = '{''thereâ€\x9d'', ''itâ€™s'', ''â€œoh'', ''cannot''}'
$ws.Cells.Item(27, 2).Value = '{''steepstaff'', ''clientele;'', ''plateall'', ''reasonablethe''}'
$ws.Cells.Item(28, 2).Value = '{''laurenour'', ''gardenoh'', ''vesuviusbreakfast'', ''valuethe'', ''pathwaysif''}'
$ws.Cells.Item(30, 2).Value = '{''kindnessthe'', '':)all'', ''absoulutely''}'
$ws.Cells.Item(32, 2).Value = '{''enderlin'', ''themour'', ''staffthe'', ''(grazie)'', ''amintas'', ''michaeli'', ''goodbyein'', ''aminta'', ''sorrentoprior'', ''emailed'', ''waitwe'', ''perfectionlucai'', ''me)'', ''futurethe'', ''standardthe'', ''seawe'', ''trainwhen'', ''hote'', ''amintasuzanne'', ''viewsalthough'', ''victorria''}'
$ws.Cells.Item(34, 2).Value = '{''issues:the'', ''wardrobebreakfast'', ''spectacularwe''}'
$ws.Cells.Item(36, 2).Value = '{''toowe'', ''didnt'', ''either)'', ''(and'', ''werent''}'
$ws.Cells.Item(37, 2).Value = '{''complexin'', ''iâ€™d'', ''motherâ€™s'', ''tripthe'', ''townthe'', ''(mind'', ''menu)'', ''aminta''}'
$ws.Cells.Item(38, 2).Value = '{''travelers)'', ''pasquales'', ''aminta'', ''(and'', ''14th''}'
$ws.Cells.Item(39, 2).Value = '{''40th'', ''towni'', ''hoteli'', ''3:30am'', ''curtiousit''}'
$ws.Cells.Item(43, 2).Value = '{''coffeethe'', ''15th'', ''couldnâ€™t''}'
$ws.Cells.Item(44, 2).Value = '{''overall--staff'', ''8:30'', ''(off-season)'', ''in-season'', ''aminta'', ''9:30'', ''helful'', ''5-10'', ''traveled''}'
$ws.Cells.Item(47, 2).Value = '{''yearthe'', ''couldnâ€™t'', ''12:30-3:00'', ''didnâ€™t'', ''daughter)'', ''weather)'', ''(as'', ''(me'', ''weâ€™d''}'
$ws.Cells.Item(48, 2).Value = '{''driverwe'', ''didnâ€™t'', ''hillreally'', ''sunset/sea'', ''aminta'', ''outside/pool''}'
$ws.Cells.Item(49, 2).Value = '{''tooi'', '':)'', ''couldnâ€™t'', ''receptionwe'', ''canâ€™t''}'
$ws.Cells.Item(51, 2).Value = '{''barstaff'', ''viewsndiwn'', ''excellemt''}'
$ws.Cells.Item(53, 2).Value = '{''before)'', ''reading)'', ''(this'', ''(you'', ''pastry/bread'', ''aminta'', ''viewsthey'', ''(especially'', ''be)'', ''nightbreakfast'', ''tripadvisor'', ''(reserved'', ''menu)'', ''(where'', ''down)'', ''course)sorrento'', ''owner(s)reading'', ''(apart'', ''dayyou''}'
$ws.Cells.Item(54, 2).Value = '{''donâ€™t'', ''impeccablythe'', ''afterit'', ''knowledgeablewe'', ''didnâ€™t'', ''hotelthere'', ''tippingwe'', ''backbreakfast'', ''poolall'', ''thinkwe''}'
$ws.Cells.Item(55, 2).Value = '{''(remember'', ''folks)'', ''barstaff/lunch'', ''â‚¬50'', ''4-5''}'
$ws.Cells.Item(62, 2).Value = '{''pastriestheres'', ''sorento'', ''overnight)'', ''didnt'', ''9:30am'', ''lantica'', ''(we'', ''labate'', ''aminta'', ''(but'', ''alleys)'', ''meal)'', ''tomato)'', ''special)we'', ''mornings;'', ''(although'', ''(steeply'', ''bad;'', ''youre'', ''onein'', ''pleasantthe'', ''favorably'', ''(though'', ''pressure;'', ''5-10''}'
$ws.Cells.Item(64, 2).Value = '{''aminta'', ''thereafterpool'', ''lotthere''}'
$ws.Cells.Item(66, 2).Value = '{''hill:'', ''rooms;'', ''staff;''}'
$ws.Cells.Item(69, 2).Value = '{''wifes'', ''prosecco'', '':)''}'
$ws.Cells.Item(71, 2).Value = '{''september?'', ''"we'', ''operatorssuch'', ''why?such'', ''plastic"'', ''1970s'', ''wordeverything'', ''dacampo'', ''slimethere'', ''tomatoesevery''}'
$ws.Cells.Item(72, 2).Value = '{''sea/mountain/garden'', ''hotelcom'', ''â‚¬174/per'', ''couldnâ€™t''}'
$ws.Cells.Item(74, 2).Value = '{''donâ€™t'', ''aminta'', ''wifi'', ''iâ€™m'', ''youâ€™d'', ''startrooms'', ''clinicalstill'', ''advisorwarmly'', ''locationsonly'', ''drinkshotel'', ''(big)'', ''(for'', ''honeymoon)'', ''youâ€™ll'', ''topoverall'', ''canâ€™t'', ''itrestaurant'', ''itâ€™s'', ''thereâ€™s''}'
$ws.Cells.Item(77, 2).Value = '{''complainthere'', ''didnt'', ''b&b'', ''poolwe'', ''atmospherethere''}'
$ws.Cells.Item(78, 2).Value = '{''900am'', ''doorbreakfast'', ''nightthe'', ''â£250'', ''soundproofing'', ''butin'', ''that;''}'
$ws.Cells.Item(81, 2).Value = '{''dessert/'', ''timethere'', ''interestingi'', ''choices-'', ''timei'', ''swimthe''}'
$ws.Cells.Item(82, 2).Value = '{''70euros'', ''itâ€™s'', ''iâ€™m''}'
$ws.Cells.Item(83, 2).Value = '{''was/offered'', ''arenâ€™t'', ''wouldnâ€™t'', ''couldnâ€™t'', ''werenâ€™t'', ''itâ€™s''}'
$ws.Cells.Item(84, 2).Value = '{''youre'', ''(and'', ''been)'', ''isnt''}'
$ws.Cells.Item(85, 2).Value = '{''showerbreakfast'', ''didnt'', ''"sunset'', ''330pm'', ''anywayvery'', ''regiona'', ''aminta'', ''view"'', ''20-25min'', ''desiredall''}'
$ws.Cells.Item(87, 2).Value = '{''sorento'', ''didnt'', ''riveria'', ''aminta'', ''traveled'', ''busthat''}'
$ws.Cells.Item(91, 2).Value = '{''doesnt'', ''inthis'', ''tiledbreakfasts'', ''cabfinal'', ''20-'', ''positivesit'', ''10:00am'', ''11:30'', ''unnecesary'', ''didnt'', ''eu150'', ''supplementbest'', ''aviod'', ''9:30'', ''available?on'', ''hankerchief''}'
$ws.Cells.Item(94, 2).Value = '{''aminta'', ''onthe'', ''itâ€™s'', ''hotelâ€™s''}'
$ws.Cells.Item(95, 2).Value = '{''viewsthe'', ''bagsroom'', ''milk)'', ''betterfood''}'
$ws.Cells.Item(97, 2).Value = '{''areanot'', ''worldcharm'', ''good:'', ''kidsthe'', ''back?'', ''points:'', ''10:30'', ''mediterranean?'', ''colleaguegood''}'
$ws.Cells.Item(100, 2).Value = '{''aminta'', ''attentivewe'', ''cannot''}'
